# Include Emissions from Imported Electricity
# Flip the boolean control lever on the BIEfIE sheet from 0 to 1.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BIEfIE")

# Set the boolean lever value (B2) to 1 -> include emissions from imported electricity
$ws.Range("B2").Value = 1

# Make BIEfIE the active sheet and select B3, matching the saved UI state
$ws.Activate()
$ws.Range("B3").Select()
